# Atualizado por script em 01-12-2023 14:45
# Re-applies the upstream scraper's re-sort of same-date fixtures and appends
# the two newly scraped matches (rows 98-99) that were added at the end of
# the sheet. Columns A:E (Indice/pais/torneio/temporada/data_partida) are
# left untouched for every existing row; only the match-detail columns
# F:V are rewritten where the diff shows a change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder matches within same-date groups (rows re-sorted by upstream scraper) ---
# row 3 <- old row 4
$ws.Range("F3").Value = "Solin"
$ws.Range("G3").Value = 0.0
$ws.Range("H3").Value = "Bijelo Brdo"
$ws.Range("I3").Value = 0.0
$ws.Range("J3").Value = 2.09
$ws.Range("K3").Value = "11/08/2023 05:43"
$ws.Range("L3").Value = 2.1
$ws.Range("M3").Value = "12/08/2023 17:27"
$ws.Range("N3").Value = 3.27
$ws.Range("O3").Value = "11/08/2023 05:43"
$ws.Range("P3").Value = 3.53
$ws.Range("Q3").Value = "12/08/2023 17:27"
$ws.Range("R3").Value = 3.12
$ws.Range("S3").Value = "11/08/2023 05:43"
$ws.Range("T3").Value = 3.2
$ws.Range("U3").Value = "12/08/2023 17:27"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/croatia/prva-nl/solin-bijelo-brdo/0vop4Qu6/"

# row 4 <- old row 5
$ws.Range("F4").Value = "Sesvete"
$ws.Range("G4").Value = 0.0
$ws.Range("H4").Value = "Dugopolje"
$ws.Range("I4").Value = 1.0
$ws.Range("J4").Value = 2.49
$ws.Range("K4").Value = "12/08/2023 17:11"
$ws.Range("L4").Value = 2.26
$ws.Range("M4").Value = "12/08/2023 17:25"
$ws.Range("N4").Value = 2.94
$ws.Range("O4").Value = "12/08/2023 17:11"
$ws.Range("P4").Value = 3.14
$ws.Range("Q4").Value = "12/08/2023 17:25"
$ws.Range("R4").Value = 2.65
$ws.Range("S4").Value = "12/08/2023 17:11"
$ws.Range("T4").Value = 3.23
$ws.Range("U4").Value = "12/08/2023 17:25"
$ws.Range("V4").Value = "https://www.betexplorer.com/football/croatia/prva-nl/sesvete-dugopolje/jNS5LRAt/"

# row 5 <- old row 6
$ws.Range("F5").Value = "Orijent"
$ws.Range("G5").Value = 1.0
$ws.Range("H5").Value = "Dubrava"
$ws.Range("I5").Value = 1.0
$ws.Range("J5").Value = 2.12
$ws.Range("K5").Value = "11/08/2023 05:43"
$ws.Range("L5").Value = 2.01
$ws.Range("M5").Value = "12/08/2023 17:29"
$ws.Range("N5").Value = 3.2
$ws.Range("O5").Value = "11/08/2023 05:43"
$ws.Range("P5").Value = 3.38
$ws.Range("Q5").Value = "12/08/2023 17:29"
$ws.Range("R5").Value = 3.12
$ws.Range("S5").Value = "11/08/2023 05:43"
$ws.Range("T5").Value = 3.59
$ws.Range("U5").Value = "12/08/2023 17:29"
$ws.Range("V5").Value = "https://www.betexplorer.com/football/croatia/prva-nl/orijent-dubrava-zagreb/Cdce1OPO/"

# row 6 <- old row 3
$ws.Range("F6").Value = "Cibalia"
$ws.Range("G6").Value = 1.0
$ws.Range("H6").Value = "Croatia Zmijavci"
$ws.Range("I6").Value = 0.0
$ws.Range("J6").Value = 1.65
$ws.Range("K6").Value = "11/08/2023 05:43"
$ws.Range("L6").Value = 2.02
$ws.Range("M6").Value = "12/08/2023 15:37"
$ws.Range("N6").Value = 3.47
$ws.Range("O6").Value = "11/08/2023 05:43"
$ws.Range("P6").Value = 3.47
$ws.Range("Q6").Value = "12/08/2023 17:16"
$ws.Range("R6").Value = 4.52
$ws.Range("S6").Value = "11/08/2023 05:43"
$ws.Range("T6").Value = 3.45
$ws.Range("U6").Value = "12/08/2023 17:16"
$ws.Range("V6").Value = "https://www.betexplorer.com/football/croatia/prva-nl/cibalia-croatia-zmijavci/E1yvQ72P/"

# row 21 <- old row 24
$ws.Range("F21").Value = "Bijelo Brdo"
$ws.Range("G21").Value = 0.0
$ws.Range("H21").Value = "Croatia Zmijavci"
$ws.Range("I21").Value = 0.0
$ws.Range("J21").Value = 2.03
$ws.Range("K21").Value = "01/09/2023 04:43"
$ws.Range("L21").Value = 2.42
$ws.Range("M21").Value = "02/09/2023 16:28"
$ws.Range("N21").Value = 3.16
$ws.Range("O21").Value = "01/09/2023 04:43"
$ws.Range("P21").Value = 3.18
$ws.Range("Q21").Value = "02/09/2023 16:28"
$ws.Range("R21").Value = 3.36
$ws.Range("S21").Value = "01/09/2023 04:43"
$ws.Range("T21").Value = 2.91
$ws.Range("U21").Value = "02/09/2023 16:28"
$ws.Range("V21").Value = "https://www.betexplorer.com/football/croatia/prva-nl/bijelo-brdo-croatia-zmijavci/rgN98fC9/"

# row 22 <- old row 21
$ws.Range("F22").Value = "Vukovar 1991"
$ws.Range("G22").Value = 3.0
$ws.Range("H22").Value = "Cibalia"
$ws.Range("I22").Value = 1.0
$ws.Range("J22").Value = 2.07
$ws.Range("K22").Value = "01/09/2023 04:43"
$ws.Range("L22").Value = 1.98
$ws.Range("M22").Value = "02/09/2023 16:28"
$ws.Range("N22").Value = 3.14
$ws.Range("O22").Value = "01/09/2023 04:43"
$ws.Range("P22").Value = 3.32
$ws.Range("Q22").Value = "02/09/2023 16:28"
$ws.Range("R22").Value = 3.3
$ws.Range("S22").Value = "01/09/2023 04:43"
$ws.Range("T22").Value = 3.76
$ws.Range("U22").Value = "02/09/2023 16:28"
$ws.Range("V22").Value = "https://www.betexplorer.com/football/croatia/prva-nl/vukovar-1991-cibalia/Srlok6Bl/"

# row 24 <- old row 22
$ws.Range("F24").Value = "Solin"
$ws.Range("G24").Value = 4.0
$ws.Range("H24").Value = "Orijent"
$ws.Range("I24").Value = 1.0
$ws.Range("J24").Value = 2.08
$ws.Range("K24").Value = "01/09/2023 04:43"
$ws.Range("L24").Value = 2.05
$ws.Range("M24").Value = "02/09/2023 16:28"
$ws.Range("N24").Value = 3.34
$ws.Range("O24").Value = "01/09/2023 04:43"
$ws.Range("P24").Value = 3.65
$ws.Range("Q24").Value = "02/09/2023 16:28"
$ws.Range("R24").Value = 2.99
$ws.Range("S24").Value = "01/09/2023 04:43"
$ws.Range("T24").Value = 3.24
$ws.Range("U24").Value = "02/09/2023 16:28"
$ws.Range("V24").Value = "https://www.betexplorer.com/football/croatia/prva-nl/solin-orijent/A7wjlQQf/"

# row 35 <- old row 36
$ws.Range("F35").Value = "Bijelo Brdo"
$ws.Range("G35").Value = 1.0
$ws.Range("H35").Value = "Dubrava"
$ws.Range("I35").Value = 1.0
$ws.Range("J35").Value = 2.17
$ws.Range("K35").Value = "15/09/2023 04:42"
$ws.Range("L35").Value = 2.54
$ws.Range("M35").Value = "16/09/2023 16:29"
$ws.Range("N35").Value = 3.12
$ws.Range("O35").Value = "15/09/2023 04:42"
$ws.Range("P35").Value = 3.11
$ws.Range("Q35").Value = "16/09/2023 16:29"
$ws.Range("R35").Value = 3.0
$ws.Range("S35").Value = "15/09/2023 04:42"
$ws.Range("T35").Value = 2.81
$ws.Range("U35").Value = "16/09/2023 16:29"
$ws.Range("V35").Value = "https://www.betexplorer.com/football/croatia/prva-nl/bijelo-brdo-dubrava-zagreb/zPuH2Pul/"

# row 36 <- old row 35
$ws.Range("F36").Value = "Solin"
$ws.Range("G36").Value = 0.0
$ws.Range("H36").Value = "Sesvete"
$ws.Range("I36").Value = 1.0
$ws.Range("J36").Value = 2.0
$ws.Range("K36").Value = "15/09/2023 04:42"
$ws.Range("L36").Value = 1.79
$ws.Range("M36").Value = "16/09/2023 16:26"
$ws.Range("N36").Value = 3.3
$ws.Range("O36").Value = "15/09/2023 04:42"
$ws.Range("P36").Value = 3.75
$ws.Range("Q36").Value = "16/09/2023 16:26"
$ws.Range("R36").Value = 3.17
$ws.Range("S36").Value = "15/09/2023 04:42"
$ws.Range("T36").Value = 4.08
$ws.Range("U36").Value = "16/09/2023 16:26"
$ws.Range("V36").Value = "https://www.betexplorer.com/football/croatia/prva-nl/solin-sesvete/MurTaNP6/"

# row 57 <- old row 59
$ws.Range("F57").Value = "Cibalia"
$ws.Range("G57").Value = 0.0
$ws.Range("H57").Value = "Sesvete"
$ws.Range("I57").Value = 1.0
$ws.Range("J57").Value = 1.73
$ws.Range("K57").Value = "13/10/2023 02:13"
$ws.Range("L57").Value = 1.75
$ws.Range("M57").Value = "14/10/2023 14:59"
$ws.Range("N57").Value = 3.53
$ws.Range("O57").Value = "13/10/2023 02:13"
$ws.Range("P57").Value = 3.52
$ws.Range("Q57").Value = "14/10/2023 14:59"
$ws.Range("R57").Value = 3.94
$ws.Range("S57").Value = "13/10/2023 02:13"
$ws.Range("T57").Value = 4.64
$ws.Range("U57").Value = "14/10/2023 14:59"
$ws.Range("V57").Value = "https://www.betexplorer.com/football/croatia/prva-nl/cibalia-sesvete/h2qHP74b/"

# row 58 <- old row 57
$ws.Range("F58").Value = "Solin"
$ws.Range("G58").Value = 1.0
$ws.Range("H58").Value = "Jarun"
$ws.Range("I58").Value = 1.0
$ws.Range("J58").Value = 1.85
$ws.Range("K58").Value = "13/10/2023 02:13"
$ws.Range("L58").Value = 1.88
$ws.Range("M58").Value = "14/10/2023 14:51"
$ws.Range("N58").Value = 3.54
$ws.Range("O58").Value = "13/10/2023 02:13"
$ws.Range("P58").Value = 3.62
$ws.Range("Q58").Value = "14/10/2023 14:51"
$ws.Range("R58").Value = 3.41
$ws.Range("S58").Value = "13/10/2023 02:13"
$ws.Range("T58").Value = 3.79
$ws.Range("U58").Value = "14/10/2023 14:51"
$ws.Range("V58").Value = "https://www.betexplorer.com/football/croatia/prva-nl/solin-jarun/6BAb7QlU/"

# row 59 <- old row 58
$ws.Range("F59").Value = "Zrinski Jurjevac"
$ws.Range("G59").Value = 4.0
$ws.Range("H59").Value = "Croatia Zmijavci"
$ws.Range("I59").Value = 0.0
$ws.Range("J59").Value = 1.68
$ws.Range("K59").Value = "13/10/2023 02:13"
$ws.Range("L59").Value = 1.58
$ws.Range("M59").Value = "14/10/2023 14:53"
$ws.Range("N59").Value = 3.61
$ws.Range("O59").Value = "13/10/2023 02:13"
$ws.Range("P59").Value = 3.86
$ws.Range("Q59").Value = "14/10/2023 14:53"
$ws.Range("R59").Value = 4.26
$ws.Range("S59").Value = "13/10/2023 02:13"
$ws.Range("T59").Value = 5.57
$ws.Range("U59").Value = "14/10/2023 14:53"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/croatia/prva-nl/zrinski-jurjevac-croatia-zmijavci/dzj8RTZo/"

# row 70 <- old row 71
$ws.Range("F70").Value = "Bijelo Brdo"
$ws.Range("G70").Value = 1.0
$ws.Range("H70").Value = "Solin"
$ws.Range("I70").Value = 1.0
$ws.Range("J70").Value = 2.14
$ws.Range("K70").Value = "27/10/2023 03:12"
$ws.Range("L70").Value = 2.68
$ws.Range("M70").Value = "28/10/2023 14:57"
$ws.Range("N70").Value = 3.17
$ws.Range("O70").Value = "27/10/2023 03:12"
$ws.Range("P70").Value = 2.66
$ws.Range("Q70").Value = "28/10/2023 14:56"
$ws.Range("R70").Value = 3.02
$ws.Range("S70").Value = "27/10/2023 03:12"
$ws.Range("T70").Value = 3.1
$ws.Range("U70").Value = "28/10/2023 14:57"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/croatia/prva-nl/bijelo-brdo-solin/hM3tS0Qd/"

# row 71 <- old row 70
$ws.Range("F71").Value = "Dugopolje"
$ws.Range("G71").Value = 1.0
$ws.Range("H71").Value = "Sesvete"
$ws.Range("I71").Value = 1.0
$ws.Range("J71").Value = 1.72
$ws.Range("K71").Value = "27/10/2023 03:12"
$ws.Range("L71").Value = 1.68
$ws.Range("M71").Value = "27/10/2023 13:14"
$ws.Range("N71").Value = 3.47
$ws.Range("O71").Value = "27/10/2023 03:12"
$ws.Range("P71").Value = 3.76
$ws.Range("Q71").Value = "28/10/2023 14:26"
$ws.Range("R71").Value = 4.07
$ws.Range("S71").Value = "27/10/2023 03:12"
$ws.Range("T71").Value = 4.74
$ws.Range("U71").Value = "28/10/2023 14:26"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-sesvete/Yya1NxuS/"

# row 76 <- old row 78
$ws.Range("F76").Value = "Orijent"
$ws.Range("G76").Value = 0.0
$ws.Range("H76").Value = "Jarun"
$ws.Range("I76").Value = 1.0
$ws.Range("J76").Value = 1.79
$ws.Range("K76").Value = "03/11/2023 02:12"
$ws.Range("L76").Value = 1.79
$ws.Range("M76").Value = "04/11/2023 13:52"
$ws.Range("N76").Value = 3.62
$ws.Range("O76").Value = "03/11/2023 02:12"
$ws.Range("P76").Value = 4.01
$ws.Range("Q76").Value = "04/11/2023 13:52"
$ws.Range("R76").Value = 3.55
$ws.Range("S76").Value = "03/11/2023 02:12"
$ws.Range("T76").Value = 3.79
$ws.Range("U76").Value = "04/11/2023 13:52"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/croatia/prva-nl/orijent-jarun/QZwtGfX1/"

# row 77 <- old row 76
$ws.Range("F77").Value = "Cibalia"
$ws.Range("G77").Value = 1.0
$ws.Range("H77").Value = "Dubrava"
$ws.Range("I77").Value = 0.0
$ws.Range("J77").Value = 1.98
$ws.Range("K77").Value = "03/11/2023 02:12"
$ws.Range("L77").Value = 2.56
$ws.Range("M77").Value = "04/11/2023 13:59"
$ws.Range("N77").Value = 3.32
$ws.Range("O77").Value = "03/11/2023 02:12"
$ws.Range("P77").Value = 3.28
$ws.Range("Q77").Value = "04/11/2023 13:59"
$ws.Range("R77").Value = 3.25
$ws.Range("S77").Value = "03/11/2023 02:12"
$ws.Range("T77").Value = 2.67
$ws.Range("U77").Value = "04/11/2023 13:59"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/croatia/prva-nl/cibalia-dubrava-zagreb/nquxHzIe/"

# row 78 <- old row 77
$ws.Range("F78").Value = "Dugopolje"
$ws.Range("G78").Value = 1.0
$ws.Range("H78").Value = "Bijelo Brdo"
$ws.Range("I78").Value = 1.0
$ws.Range("J78").Value = 1.62
$ws.Range("K78").Value = "03/11/2023 02:12"
$ws.Range("L78").Value = 1.92
$ws.Range("M78").Value = "04/11/2023 13:59"
$ws.Range("N78").Value = 3.57
$ws.Range("O78").Value = "03/11/2023 02:12"
$ws.Range("P78").Value = 3.03
$ws.Range("Q78").Value = "04/11/2023 13:59"
$ws.Range("R78").Value = 4.6
$ws.Range("S78").Value = "03/11/2023 02:12"
$ws.Range("T78").Value = 4.54
$ws.Range("U78").Value = "04/11/2023 13:59"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-bijelo-brdo/jHoTIdmq/"

# row 86 <- old row 87
$ws.Range("F86").Value = "Sesvete"
$ws.Range("G86").Value = 2.0
$ws.Range("H86").Value = "Jarun"
$ws.Range("I86").Value = 1.0
$ws.Range("J86").Value = 1.76
$ws.Range("K86").Value = "16/11/2023 01:42"
$ws.Range("L86").Value = 1.79
$ws.Range("M86").Value = "17/11/2023 13:29"
$ws.Range("N86").Value = 3.49
$ws.Range("O86").Value = "16/11/2023 01:42"
$ws.Range("P86").Value = 3.09
$ws.Range("Q86").Value = "17/11/2023 13:29"
$ws.Range("R86").Value = 3.84
$ws.Range("S86").Value = "16/11/2023 01:42"
$ws.Range("T86").Value = 5.3
$ws.Range("U86").Value = "17/11/2023 13:29"
$ws.Range("V86").Value = "https://www.betexplorer.com/football/croatia/prva-nl/sesvete-jarun/Gxd5iHg8/"

# row 87 <- old row 86
$ws.Range("F87").Value = "Cibalia"
$ws.Range("G87").Value = 0.0
$ws.Range("H87").Value = "Vukovar 1991"
$ws.Range("I87").Value = 1.0
$ws.Range("J87").Value = 2.61
$ws.Range("K87").Value = "16/11/2023 01:42"
$ws.Range("L87").Value = 3.57
$ws.Range("M87").Value = "17/11/2023 13:25"
$ws.Range("N87").Value = 3.09
$ws.Range("O87").Value = "16/11/2023 01:42"
$ws.Range("P87").Value = 3.31
$ws.Range("Q87").Value = "17/11/2023 13:25"
$ws.Range("R87").Value = 2.52
$ws.Range("S87").Value = "16/11/2023 01:42"
$ws.Range("T87").Value = 2.04
$ws.Range("U87").Value = "17/11/2023 13:25"
$ws.Range("V87").Value = "https://www.betexplorer.com/football/croatia/prva-nl/cibalia-vukovar-1991/fTf9jy9E/"

# row 88 <- old row 89
$ws.Range("F88").Value = "Orijent"
$ws.Range("G88").Value = 1.0
$ws.Range("H88").Value = "Solin"
$ws.Range("I88").Value = 1.0
$ws.Range("J88").Value = 1.95
$ws.Range("K88").Value = "17/11/2023 01:42"
$ws.Range("L88").Value = 1.88
$ws.Range("M88").Value = "18/11/2023 13:25"
$ws.Range("N88").Value = 3.43
$ws.Range("O88").Value = "17/11/2023 01:42"
$ws.Range("P88").Value = 3.71
$ws.Range("Q88").Value = "18/11/2023 13:25"
$ws.Range("R88").Value = 3.22
$ws.Range("S88").Value = "17/11/2023 01:42"
$ws.Range("T88").Value = 3.67
$ws.Range("U88").Value = "18/11/2023 13:25"
$ws.Range("V88").Value = "https://www.betexplorer.com/football/croatia/prva-nl/orijent-solin/ljfDkeOK/"

# row 89 <- old row 88
$ws.Range("F89").Value = "Dugopolje"
$ws.Range("G89").Value = 1.0
$ws.Range("H89").Value = "Dubrava"
$ws.Range("I89").Value = 2.0
$ws.Range("J89").Value = 1.81
$ws.Range("K89").Value = "17/11/2023 01:42"
$ws.Range("L89").Value = 2.05
$ws.Range("M89").Value = "18/11/2023 13:21"
$ws.Range("N89").Value = 3.44
$ws.Range("O89").Value = "17/11/2023 01:42"
$ws.Range("P89").Value = 3.39
$ws.Range("Q89").Value = "18/11/2023 13:21"
$ws.Range("R89").Value = 3.81
$ws.Range("S89").Value = "17/11/2023 01:42"
$ws.Range("T89").Value = 3.48
$ws.Range("U89").Value = "18/11/2023 13:21"
$ws.Range("V89").Value = "https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-dubrava-zagreb/vLCMSXG7/"

# --- Append new rows 98 and 99 (new matches added by scraper) ---
$ws.Range("A97:V97").Copy($ws.Range("A98:V99"))

# row 98
$ws.Range("A98").Value = 97.0
$ws.Range("B98").Value = "croatia"
$ws.Range("C98").Value = "prva-nl"
$ws.Range("D98").Value = "2023-2024"
$ws.Range("E98").Value = 45261.54166666666
$ws.Range("F98").Value = "Dubrava"
$ws.Range("G98").Value = 2.0
$ws.Range("H98").Value = "Bijelo Brdo"
$ws.Range("I98").Value = 1.0
$ws.Range("J98").Value = 1.79
$ws.Range("K98").Value = "30/11/2023 01:12"
$ws.Range("L98").Value = 1.74
$ws.Range("M98").Value = "01/12/2023 12:58"
$ws.Range("N98").Value = 3.33
$ws.Range("O98").Value = "30/11/2023 01:12"
$ws.Range("P98").Value = 3.52
$ws.Range("Q98").Value = "01/12/2023 12:58"
$ws.Range("R98").Value = 4.0
$ws.Range("S98").Value = "30/11/2023 01:12"
$ws.Range("T98").Value = 4.75
$ws.Range("U98").Value = "01/12/2023 12:58"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/croatia/prva-nl/dubrava-zagreb-bijelo-brdo/tpJj3EhK/"

# row 99
$ws.Range("A99").Value = 98.0
$ws.Range("B99").Value = "croatia"
$ws.Range("C99").Value = "prva-nl"
$ws.Range("D99").Value = "2023-2024"
$ws.Range("E99").Value = 45261.5625
$ws.Range("F99").Value = "Sesvete"
$ws.Range("G99").Value = 2.0
$ws.Range("H99").Value = "Solin"
$ws.Range("I99").Value = 0.0
$ws.Range("J99").Value = 1.96
$ws.Range("K99").Value = "30/11/2023 01:42"
$ws.Range("L99").Value = 1.99
$ws.Range("M99").Value = "01/12/2023 13:22"
$ws.Range("N99").Value = 3.25
$ws.Range("O99").Value = "30/11/2023 01:42"
$ws.Range("P99").Value = 3.5
$ws.Range("Q99").Value = "01/12/2023 13:05"
$ws.Range("R99").Value = 3.45
$ws.Range("S99").Value = "30/11/2023 01:42"
$ws.Range("T99").Value = 3.55
$ws.Range("U99").Value = "01/12/2023 13:22"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/croatia/prva-nl/sesvete-solin/IRLJdCVm/"
